$p = $ppt.ActivePresentation

function Set-RunBold {
    param($shape, [string]$needle, [int]$bold)
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($needle)
    if ($idx -lt 0) { return }
    $sub = $tr.Characters($idx + 1, $needle.Length)
    $sub.Font.Bold = $bold
}

function Replace-Text {
    param($shape, [string]$needle, [string]$replacement)
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($needle)
    if ($idx -lt 0) { return }
    $sub = $tr.Characters($idx + 1, $needle.Length)
    $sub.Text = $replacement
}

# ---------- Slide 3: "running time of 401.7s" -> split & extend ----------
$s3 = $p.Slides.Item(3)
$grp3 = $s3.Shapes.Item(4)
$rect7 = $grp3.GroupItems.Item(3)
Replace-Text $rect7 "running time of 401.7s" "running time of 401.7s on the test set."
Set-RunBold $rect7 "on the test set." 0

# ---------- Slide 5: HSV "ze" -> "we" typo fix ----------
$s5 = $p.Slides.Item(5)
$hsvBox = $s5.Shapes.Item(7)
Replace-Text $hsvBox "ze used" "we used"

# ---------- Slide 7: position adjustments ----------
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(2).Left = 118529 / 914400 * 72
$s7.Shapes.Item(4).Left = 118529 / 914400 * 72
$s7.Shapes.Item(3).Left = 514350 / 914400 * 72
$s7.Shapes.Item(5).Left = 514350 / 914400 * 72
$s7.Shapes.Item(5).Width = 6800850 / 914400 * 72

# ---------- Slide 7: text adjustments ----------
$xgbBox = $s7.Shapes.Item(5)
Replace-Text $xgbBox "we used Cross-Validation to fine the best parameters: " "we used Cross-Validation on the training set to fine the best parameters: "
Replace-Text $xgbBox "we select  " "we select ed "
Replace-Text $xgbBox "The running time is 56s." "The running time is 56son the test set."
